# Sold Dark Petition to a hot guy
# Remove row 36 (Dark Petition / Magic Origins) entirely; remaining rows shift up by one.
# The SUM formula and sheet dimension auto-adjust when the row is deleted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:36").Delete()

# Update card prices (column D) for the remaining 98 rows to their latest values.
$ws.Range("D2").Value = 9.76
$ws.Range("D3").Value = 14.25
$ws.Range("D4").Value = 17.13
$ws.Range("D5").Value = 9.26
$ws.Range("D6").Value = 8.77
$ws.Range("D7").Value = 0.59
$ws.Range("D8").Value = 0.95
$ws.Range("D9").Value = 1.63
$ws.Range("D10").Value = 1.19
$ws.Range("D11").Value = 3.38
$ws.Range("D12").Value = 2.65
$ws.Range("D13").Value = 0.77
$ws.Range("D14").Value = 0.24
$ws.Range("D15").Value = 0.73
$ws.Range("D16").Value = 0.1
$ws.Range("D17").Value = 0.1
$ws.Range("D18").Value = 0.59
$ws.Range("D19").Value = 7.81
$ws.Range("D20").Value = 2.87
$ws.Range("D21").Value = 2.99
$ws.Range("D22").Value = 0.75
$ws.Range("D23").Value = 3.64
$ws.Range("D24").Value = 30.54
$ws.Range("D25").Value = 5.8
$ws.Range("D26").Value = 20.05
$ws.Range("D27").Value = 2.13
$ws.Range("D28").Value = 0.64
$ws.Range("D29").Value = 4.41
$ws.Range("D30").Value = 4.07
$ws.Range("D31").Value = 0.75
$ws.Range("D32").Value = 0.6899999999999999
$ws.Range("D33").Value = 0.83
$ws.Range("D34").Value = 0.98
$ws.Range("D35").Value = 15.1
$ws.Range("D36").Value = 6.61
$ws.Range("D37").Value = 8.08
$ws.Range("D38").Value = 4.14
$ws.Range("D39").Value = 0.67
$ws.Range("D40").Value = 1.33
$ws.Range("D41").Value = 54.62
$ws.Range("D42").Value = 0.96
$ws.Range("D43").Value = 0.96
$ws.Range("D44").Value = 5.06
$ws.Range("D45").Value = 0.99
$ws.Range("D46").Value = 3.39
$ws.Range("D47").Value = 6.07
$ws.Range("D48").Value = 1.49
$ws.Range("D49").Value = 1.6
$ws.Range("D50").Value = 4.95
$ws.Range("D51").Value = 0.87
$ws.Range("D52").Value = 4.17
$ws.Range("D53").Value = 3.14
$ws.Range("D54").Value = 1.56
$ws.Range("D55").Value = 1.3
$ws.Range("D56").Value = 1.28
$ws.Range("D57").Value = 1.28
$ws.Range("D58").Value = 1.61
$ws.Range("D59").Value = 1.73
$ws.Range("D60").Value = 2.63
$ws.Range("D61").Value = 3.28
$ws.Range("D62").Value = 1.83
$ws.Range("D63").Value = 4.95
$ws.Range("D64").Value = 27.63
$ws.Range("D65").Value = 9.49
$ws.Range("D66").Value = 1.56
$ws.Range("D67").Value = 1
$ws.Range("D68").Value = 3.49
$ws.Range("D69").Value = 4.19
$ws.Range("D70").Value = 0.64
$ws.Range("D71").Value = 0.76
$ws.Range("D72").Value = 2.6
$ws.Range("D73").Value = 2.95
$ws.Range("D74").Value = 6.76
$ws.Range("D75").Value = 4.94
$ws.Range("D76").Value = 0.96
$ws.Range("D77").Value = 1.21
$ws.Range("D78").Value = 6.07
$ws.Range("D79").Value = 3.76
$ws.Range("D80").Value = 3.47
$ws.Range("D81").Value = 1.44
$ws.Range("D82").Value = 1.38
$ws.Range("D83").Value = 9.6
$ws.Range("D84").Value = 7.28
$ws.Range("D85").Value = 18.29
$ws.Range("D86").Value = 1.07
$ws.Range("D87").Value = 9.699999999999999
$ws.Range("D88").Value = 9.33
$ws.Range("D89").Value = 3.28
$ws.Range("D90").Value = 1.9
$ws.Range("D91").Value = 0.61
$ws.Range("D92").Value = 4.82
$ws.Range("D93").Value = 3.61
$ws.Range("D94").Value = 0.75
$ws.Range("D95").Value = 0.59
$ws.Range("D96").Value = 1.76
$ws.Range("D97").Value = 0.85
$ws.Range("D98").Value = 0.6899999999999999
$ws.Range("D99").Value = 1.36

# Update the sheet view to reflect scrolling to around row 25 with row 36 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows("36:36").Select()
